# Update Excel file with latest predictions

$wb = $excel.ActiveWorkbook

# --- Sheet "Home win" ---
$ws1 = $wb.Worksheets.Item("Home win")
$ws1.Range("A2").Value = "24-12-2024 10:00"
$ws1.Range("B2").Value = "AZERBAIDJAN"
$ws1.Range("C2").Value = "BIRINCI DASTA"
$ws1.Range("D2").Value = "Qaradağ Lökbatan - Zaqatala"
$ws1.Range("F2").Value = 3.6

# --- Sheet "Draw" ---
$ws2 = $wb.Worksheets.Item("Draw")
$ws2.Range("A2").Value = "25-12-2024 17:30"
$ws2.Range("B2").Value = "WORLD"
$ws2.Range("C2").Value = "GULF CUP OF NATIONS"
$ws2.Range("D2").Value = "Bahrain - Iraq"
$ws2.Range("E2").Value = 60
$ws2.Range("F2").Value = 2.9

# --- Sheet "Btts" ---
$ws3 = $wb.Worksheets.Item("Btts")
$ws3.Range("A2").Value = "24-12-2024 17:00"
$ws3.Range("B2").Value = "ISRAEL"
$ws3.Range("C2").Value = "LIGA ALEF"
$ws3.Range("D2").Value = "Tira - Maccabi Ahi Nazareth"
$ws3.Range("E2").Value = 83.3
$ws3.Range("F2").Value = 1.9
